# Apply updated dSF (column F) values to reflect repulled data / push all
# data / mean calculation changes described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = -3
    "F5"  = 1
    "F8"  = -4
    "F10" = -6
    "F12" = -9
    "F14" = -6
    "F15" = 3
    "F19" = -1
    "F20" = -4
    "F29" = -1
    "F32" = -4
    "F34" = 0
    "F38" = -1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
